$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the sheet with 11 more blank, formatted rows (11-21), reusing the
# formatting of the existing blank row 10 (same style index, no new styles
# get created in xl/styles.xml).
$ws.Range("A10:G10").Copy()
$ws.Range("A11:G21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New "Datenquelle" / source row values (row 11)
$ws.Range("A11").Value = "Datenquelle:"
$ws.Range("B11").Value = "Weltgesundheitsorganisation (WHO)"

# Drop the now-unused column H entirely (shifts nothing, just removes it)
$ws.Columns.Item(8).Delete()

# Column A/B widths to fit the new "Datenquelle:" / "Weltgesundheitsorganisation (WHO)" labels
$ws.Columns.Item(1).ColumnWidth = 11.6
$ws.Columns.Item(2).ColumnWidth = 11.6

# Match the recorded selection left behind in the source workbook
$ws.Range("J12").Select()
